$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.3304176666666667
$ws.Range("H2").Value = 0.9912529999999999

$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.3318626666666666
$ws.Range("N2").Value = 0.9955879999999999

$ws.Range("Q2").Value = 0.1096532879737778
$ws.Range("R2").Value = 0.9868795917639999
